# Add a new worksheet "APIData" as the last tab (after "Formula"), fill it
# with the symbol list used by the web/REST test suite, and make it the
# active/selected sheet - matching the target workbook.xml / sheet5.xml.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the current last sheet so it lands at the end
# of the tab strip (Worksheets.Add defaults to inserting before the active
# sheet, which is not what we want here).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "APIData"

$ws.Range("A1").Value = "symbol"
$ws.Range("A2").Value = "INFY"
$ws.Range("A3").Value = "ADSL"

# Leave the selection on the last populated cell and make this the active
# tab, mirroring the authored workbook's UI state.
$ws.Activate()
$ws.Range("A3").Select() | Out-Null
